$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.975.67'
$ws.Range("E2").Value = '  -0.89%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.561.60'
$ws.Range("E3").Value = '  -1.65%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.21'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '197.22'
$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("E7").Value = '  -2.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.210'
$ws.Range("E9").Value = '  +1.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.628'
$ws.Range("E10").Value = '  -2.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.03'
$ws.Range("E11").Value = '  -1.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000290'
$ws.Range("E12").Value = '  -4.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.29'
$ws.Range("E13").Value = '  -2.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.116.49'
$ws.Range("E14").Value = '  -1.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '659.93'
$ws.Range("E15").Value = '  +9.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.807.43'
$ws.Range("E16").Value = '  -1.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.62'
$ws.Range("E17").Value = '  -3.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.558.47'
$ws.Range("E18").Value = '  -1.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.122'
$ws.Range("E19").Value = '  -0.97%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.51'
$ws.Range("E20").Value = '  -3.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.967'
$ws.Range("E21").Value = '  -3.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.24'
$ws.Range("E22").Value = '  +2.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.37'
$ws.Range("E23").Value = '  +3.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '104.76'
$ws.Range("E24").Value = '  +2.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.41'
$ws.Range("E25").Value = '  -4.70%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.94'
$ws.Range("E26").Value = '  -3.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.27'
$ws.Range("E27").Value = '  -4.64%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.65'
$ws.Range("E28").Value = '  +0.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.46'
$ws.Range("E29").Value = '  -1.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.41'
$ws.Range("E30").Value = '  -5.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.83'
$ws.Range("E31").Value = '  -5.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.80'
$ws.Range("E32").Value = '  -4.09%  '

$ws.Range("E33").Value = '  -5.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '61.84'

$ws.Range("B35").Value = 'Stacks'
$ws.Range("C35").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.79'
$ws.Range("E35").Value = '  +6.53%  '

$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.758.69'
$ws.Range("E36").Value = '  -4.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0₃0819'
$ws.Range("E37").Value = '  -8.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '512.78'
$ws.Range("E39").Value = '  -5.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.95'
$ws.Range("E40").Value = '  -5.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.374'
$ws.Range("E41").Value = '  -4.38%  '

$ws.Range("E42").Value = '  +1.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.03'
$ws.Range("E43").Value = '  -5.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0454'
$ws.Range("E44").Value = '  -1.41%  '

$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.90'
$ws.Range("E45").Value = '  +1.19%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.41'
$ws.Range("E46").Value = '  +0.36%  '

$ws.Range("E47").Value = '  -2.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.998'
$ws.Range("E48").Value = '  -0.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.38'
$ws.Range("E49").Value = '  -2.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.75'
$ws.Range("E50").Value = '  +65.58%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.76'
$ws.Range("E51").Value = '  +17.98%  '
